$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 13 (Perspicacité): modificateur changes from Intéligence to Sagesse, bonus 3 -> 0
$ws.Range("E13").Value = "Sagesse"
$ws.Range("F13").Value = 0

# Row 15 (Investigation): modificateur changes from Sagesse to Intéligence, bonus -1 -> 3
$ws.Range("E15").Value = "Intéligence"
$ws.Range("F15").Value = 3

# Update the active cell selection to H29
$ws.Activate()
$ws.Range("H29").Select()
